# Update the "取得日時" (acquired timestamp) column for the existing
# rows of data on the "ランサーズ" sheet to reflect the latest
# scrape run at 2025-09-05 06:32:35.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-05 06:32:35"

for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
